# Add remaining SOCP constraints
# Populate the Qmax (K) / Qmin (L) columns with 0 for the storage units
# in rows 8-12 of the "Power Storage" sheet (row 7 already had values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Power Storage")

$ws.Range("K8:L12").Value = 0

# Move the active selection on the frozen (bottom-left) pane to L16,
# matching where the author left the cursor after entering the values.
[void]$ws.Range("L16").Select()
